$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.238.99'
$ws.Range("E2").Value = '  +9.19%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.240.19'
$ws.Range("E3").Value = '  +4.57%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '402.74'
$ws.Range("E5").Value = '  +4.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.42'
$ws.Range("E6").Value = '  +8.34%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.558'
$ws.Range("E7").Value = '  +3.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.627'
$ws.Range("E9").Value = '  +7.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.46'
$ws.Range("E10").Value = '  +6.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0899'
$ws.Range("E11").Value = '  +5.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.140'
$ws.Range("E12").Value = '  +2.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.759.80'
$ws.Range("E13").Value = '  +4.81%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.13'
$ws.Range("E14").Value = '  +3.72%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.17'
$ws.Range("E15").Value = '  +2.92%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.253.65'
$ws.Range("E16").Value = '  +4.90%  '

$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.07'
$ws.Range("E17").Value = '  +7.42%  '

$ws.Range("E18").Value = '  -3.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '56.225.26'
$ws.Range("E19").Value = '  +9.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.43'
$ws.Range("E20").Value = '  +2.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000103'
$ws.Range("E21").Value = '  +6.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.13'
$ws.Range("E22").Value = '  +6.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '305.52'
$ws.Range("E23").Value = '  +14.78%  '

$ws.Range("E24").Value = '  +8.13%  '

$ws.Range("E25").Value = '  +3.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.26'
$ws.Range("E26").Value = '  +1.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.42'
$ws.Range("E27").Value = '  +5.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.45'
$ws.Range("E28").Value = '  +2.85%  '

$ws.Range("E29").Value = '  +3.31%  '

$ws.Range("E31").Value = '  +5.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.41'
$ws.Range("E32").Value = '  +10.58%  '

$ws.Range("E33").Value = '  +5.86%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '36.70'
$ws.Range("E34").Value = '  +3.81%  '

$ws.Range("E35").Value = '  +3.27%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.49'
$ws.Range("E36").Value = '  +2.46%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.56'
$ws.Range("E37").Value = '  +5.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.10'
$ws.Range("E38").Value = '  +23.34%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.93'
$ws.Range("E40").Value = '  +2.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.03'
$ws.Range("E41").Value = '  +10.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.22'
$ws.Range("E42").Value = '  +3.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '131.87'
$ws.Range("E43").Value = '  +2.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.119'
$ws.Range("E44").Value = '  +3.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.285'
$ws.Range("E45").Value = '  -4.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.70'
$ws.Range("E46").Value = '  +1.19%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.162.42'
$ws.Range("E47").Value = '  +5.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.10'
$ws.Range("E48").Value = '  +41.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.10'
$ws.Range("E49").Value = '  +1.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.45'
$ws.Range("E50").Value = '  -0.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0360'
$ws.Range("E51").Value = '  +8.83%  '
